$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
